$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D2 value (typo fix: missing leading digit corrected from 46390 to 47390)
$ws.Range("D2").Value = 47390

# Clear the "Image Name" index column for rows 3-11 (missing judgement condition
# meant these rows were written with stale index values)
$ws.Range("A3:A11").ClearContents()

# Reset the measurement columns for rows 3-11 back to zero, matching the
# default/untouched rows further down the sheet
$ws.Range("B3:D11").Value = 0
